# Update "人气/浏览量" (column F) values on the "展览" and "全部类型" worksheets.
# These two sheets mirror the same underlying rows, so the same row/value
# updates are applied to both.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 2;  Value = 1067 },
    @{ Row = 3;  Value = 761 },
    @{ Row = 6;  Value = 1100 },
    @{ Row = 7;  Value = 230 },
    @{ Row = 8;  Value = 1798 },
    @{ Row = 9;  Value = 6443 },
    @{ Row = 10; Value = 490 },
    @{ Row = 11; Value = 380 },
    @{ Row = 13; Value = 108 },
    @{ Row = 14; Value = 381 },
    @{ Row = 16; Value = 6940 },
    @{ Row = 17; Value = 279 },
    @{ Row = 22; Value = 110 },
    @{ Row = 24; Value = 113 },
    @{ Row = 28; Value = 13 },
    @{ Row = 29; Value = 397 },
    @{ Row = 30; Value = 597 },
    @{ Row = 32; Value = 80 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Cells.Item($u.Row, 6).Value = $u.Value
    }
}
